$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"   = 264.000
    "H2"   = 955.15
    "G3"   = 1725.000
    "H3"   = 233235.85
    "G11"  = 3983.000
    "H11"  = 122197.09
    "G16"  = 1632.000
    "H16"  = 139483.01
    "G20"  = 2652.000
    "H20"  = 224470.60
    "H30"  = 7194.42
    "H72"  = 9934.12
    "G87"  = 4019.000
    "H87"  = 73307.88
    "G136" = 249.000
    "H136" = 83617.55
    "G140" = 191.000
    "H140" = 8972.17
    "G192" = 3455.000
    "H192" = 34733.84
    "G211" = 69534.800
    "H211" = 163009.49
    "G222" = 6081.000
    "H222" = 51446.08
    "H255" = 2313773.93
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
